$d = $word.ActiveDocument

$d.Content.Find.Execute("95×81=7695", $true, $false, $false, $false, $false, $true, 1, $false, "47×76=3572", 2) | Out-Null
$d.Content.Find.Execute("37×78=2886", $true, $false, $false, $false, $false, $true, 1, $false, "37×79=2923", 2) | Out-Null
$d.Content.Find.Execute("65×66=4290", $true, $false, $false, $false, $false, $true, 1, $false, "14×13=182", 2) | Out-Null
$d.Content.Find.Execute("46×40=1840", $true, $false, $false, $false, $false, $true, 1, $false, "61×25=1525", 2) | Out-Null
$d.Content.Find.Execute("15×98=1470", $true, $false, $false, $false, $false, $true, 1, $false, "11×19=209", 2) | Out-Null
$d.Content.Find.Execute("19×51=969", $true, $false, $false, $false, $false, $true, 1, $false, "23×51=1173", 2) | Out-Null
$d.Content.Find.Execute("58×53=3074", $true, $false, $false, $false, $false, $true, 1, $false, "18×86=1548", 2) | Out-Null
$d.Content.Find.Execute("64×59=3776", $true, $false, $false, $false, $false, $true, 1, $false, "50×42=2100", 2) | Out-Null
$d.Content.Find.Execute("74×71=5254", $true, $false, $false, $false, $false, $true, 1, $false, "13×22=286", 2) | Out-Null
$d.Content.Find.Execute("82×45=3690", $true, $false, $false, $false, $false, $true, 1, $false, "22×52=1144", 2) | Out-Null
$d.Content.Find.Execute("33×66=2178", $true, $false, $false, $false, $false, $true, 1, $false, "67×68=4556", 2) | Out-Null
$d.Content.Find.Execute("61×86=5246", $true, $false, $false, $false, $false, $true, 1, $false, "82×75=6150", 2) | Out-Null
$d.Content.Find.Execute("25×13=325", $true, $false, $false, $false, $false, $true, 1, $false, "28×56=1568", 2) | Out-Null
$d.Content.Find.Execute("80×59=4720", $true, $false, $false, $false, $false, $true, 1, $false, "92×84=7728", 2) | Out-Null
$d.Content.Find.Execute("79×44=3476", $true, $false, $false, $false, $false, $true, 1, $false, "71×99=7029", 2) | Out-Null
$d.Content.Find.Execute("80×86=6880", $true, $false, $false, $false, $false, $true, 1, $false, "63×94=5922", 2) | Out-Null
$d.Content.Find.Execute("56×52=2912", $true, $false, $false, $false, $false, $true, 1, $false, "38×43=1634", 2) | Out-Null
$d.Content.Find.Execute("77×37=2849", $true, $false, $false, $false, $false, $true, 1, $false, "28×44=1232", 2) | Out-Null
$d.Content.Find.Execute("14×92=1288", $true, $false, $false, $false, $false, $true, 1, $false, "93×45=4185", 2) | Out-Null
$d.Content.Find.Execute("87×16=1392", $true, $false, $false, $false, $false, $true, 1, $false, "49×28=1372", 2) | Out-Null
$d.Content.Find.Execute("28×20=560", $true, $false, $false, $false, $false, $true, 1, $false, "49×40=1960", 2) | Out-Null
$d.Content.Find.Execute("25×55=1375", $true, $false, $false, $false, $false, $true, 1, $false, "13×27=351", 2) | Out-Null
$d.Content.Find.Execute("19×30=570", $true, $false, $false, $false, $false, $true, 1, $false, "90×62=5580", 2) | Out-Null
$d.Content.Find.Execute("95×12=1140", $true, $false, $false, $false, $false, $true, 1, $false, "61×66=4026", 2) | Out-Null
$d.Content.Find.Execute("72×85=6120", $true, $false, $false, $false, $false, $true, 1, $false, "95×98=9310", 2) | Out-Null
